$d = $word.ActiveDocument

# --- 1) Update the 2010 National Security Strategy file-path/title line ---
$find1 = "Files\\2011 Case Study\\Primary Sources_Policy_Strategies\\2010_national_security_strategy - § 3 references coded [ 0.20% Coverage]"
$repl1 = "Files\\2011 Case Study\\CS1_Primary Sources_Policy_Strategies\\2010 National Security Strategy - § 3 references coded [ 0.20% Coverage]"
$d.Content.Find.Execute($find1, $true, $false, $false, $false, $false, $true, 1, $false, $repl1, 2)

# --- 2) Update the 2011 International Strategy for Cyberspace file-path/title line ---
$find2 = "Files\\2011 Case Study\\Primary Sources_Policy_Strategies\\2011_International_strategy_for_cyberspace - § 4 references coded [ 0.10% Coverage]"
$repl2 = "Files\\2011 Case Study\\CS1_Primary Sources_Policy_Strategies\\2011 International Strategy for Cyberspace - § 4 references coded [ 0.10% Coverage]"
$d.Content.Find.Execute($find2, $true, $false, $false, $false, $false, $true, 1, $false, $repl2, 2)

# --- 3) Update the 2015 National Security Strategy file-path/title line ---
$find3 = "Files\\2015 Case Study\\Primary Sources_Policy_Strategies\\2015 National Security Strategy CLEAN - § 9 references coded [ 0.90% Coverage]"
$repl3 = "Files\\2015 Case Study\\CS2_Primary Sources_Policy_Strategies\\2015 National Security Strategy - § 9 references coded [ 0.90% Coverage]"
$d.Content.Find.Execute($find3, $true, $false, $false, $false, $false, $true, 1, $false, $repl3, 2)

# --- 4) Append three new paragraphs (CS3 / 2017 NSS block) at the end of the document ---
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$tail = $d.Range($lastPara.Range.End, $lastPara.Range.End)

$newParasXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="TextBody"/><w:bidi w:val="0"/><w:spacing w:before="113" w:after="113"/><w:ind w:left="113" w:right="113" w:hanging="0"/><w:jc w:val="left"/><w:rPr><w:highlight w:val="lightGray"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>Files\\2018 Case Study\\CS3_Primary Sources_Policy_Strategies\\2017 National Security Strategy - § 1 reference coded [ 0.06% Coverage]</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="TextBody"/><w:bidi w:val="0"/><w:spacing w:before="113" w:after="113"/><w:ind w:left="113" w:right="113" w:hanging="0"/><w:jc w:val="left"/><w:rPr><w:highlight w:val="lightGray"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>Reference 1 - 0.06% Coverage</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="TextBody"/><w:bidi w:val="0"/><w:spacing w:before="0" w:after="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t>We will bring about the bett er future we seek for our people and the world, by confronting the challenges and dangers posed by those who seek to destabilize the world and threaten America’s people and interests.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$tail.InsertXML($newParasXml)
